$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.71
$ws.Range("C2").Value = 0.511019023580748
$ws.Range("D2").Value = 0.5
$ws.Range("G2").Value = 0.05675483310665119
$ws.Range("H2").Value = 4.226798047848368
$ws.Range("L2").Value = 93.62

# Row 3
$ws.Range("A3").Value = 489
$ws.Range("B3").Value = 22.11
$ws.Range("C3").Value = 2.879129582323678
$ws.Range("E3").Value = 22.63588915782518
$ws.Range("G3").Value = 1.761935450383476
$ws.Range("H3").Value = 66.28619713392969
$ws.Range("J3").Value = 1.380962440290202
$ws.Range("L3").Value = 319.93

# Row 4
$ws.Range("A4").Value = 481
$ws.Range("B4").Value = 53.6
$ws.Range("C4").Value = 4.519986719273934
$ws.Range("E4").Value = 36.09334023314045
$ws.Range("F4").Value = 36.03
$ws.Range("G4").Value = 2.953300250376031
$ws.Range("H4").Value = 134.4745860826723
$ws.Range("J4").Value = 2.801553876722341
$ws.Range("L4").Value = 974.9
